$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 1 (header/subject numbers) for columns B:E
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Update row 2 (CON) for columns B:E
$ws.Range("B2").Value = 0.95545643785773215
$ws.Range("C2").Value = 1.1708474476219339
$ws.Range("D2").Value = 3.4382121402409167
$ws.Range("E2").Value = 0.99020758740057668

# Update row 3 (STR) for columns B:E
$ws.Range("B3").Value = 1.5912582131690178
$ws.Range("C3").Value = 0.58061146085824544
$ws.Range("D3").Value = 2.5573480464496194
$ws.Range("E3").Value = 0.21395928295942285

# Update the selection to match the new reduced selection range
$ws.Range("B1:E3").Select()
